# Refresh the Csf2-Csf2rb LR-pair sheet with the re-run ("new TPM") NATMI
# numbers. Rows 2-11 (FAPs / MuSCs as sending cluster) get updated
# expression/specificity figures, and rows 12-16 are added for a third
# sending cluster, Neutrophils, against the same four target clusters
# (ECs, Inflammatory-Mac, MuSCs, Neutrophils, Resolving-Mac).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Csf2"
$ws.Range("C2").Value = "Csf2rb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2300786666666667
$ws.Range("H2").Value = 0.690236
$ws.Range("I2").Value = 0.2137022699341201
$ws.Range("J2").Value = 0.2304051507653011
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.595825
$ws.Range("N2").Value = 27.19165
$ws.Range("O2").Value = 0.05267224809548414
$ws.Range("P2").Value = 0.03575312334118911
$ws.Range("Q2").Value = 3.128109288233333
$ws.Range("R2").Value = 18.7686557294
$ws.Range("S2").Value = 0.0112561789805381
$ws.Range("T2").Value = 0.008237703773757082

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Csf2"
$ws.Range("C3").Value = "Csf2rb"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2300786666666667
$ws.Range("H3").Value = 0.690236
$ws.Range("I3").Value = 0.2137022699341201
$ws.Range("J3").Value = 0.2304051507653011
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 82.69899733333334
$ws.Range("N3").Value = 248.096992
$ws.Range("O3").Value = 0.3203882151167078
$ws.Range("P3").Value = 0.3262119935919301
$ws.Range("Q3").Value = 19.02727504112356
$ws.Range("R3").Value = 171.245475370112
$ws.Range("S3").Value = 0.06846768883058164
$ws.Range("T3").Value = 0.0751609235649981

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Csf2"
$ws.Range("C4").Value = "Csf2rb"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2300786666666667
$ws.Range("H4").Value = 0.690236
$ws.Range("I4").Value = 0.2137022699341201
$ws.Range("J4").Value = 0.2304051507653011
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.2286895
$ws.Range("N4").Value = 0.457379
$ws.Range("O4").Value = 0.0008859771349537243
$ws.Range("P4").Value = 0.0006013878451903335
$ws.Range("Q4").Value = 0.05261657524066666
$ws.Range("R4").Value = 0.315699451444
$ws.Range("S4").Value = 0.0001893353248493392
$ws.Range("T4").Value = 0.0001385628571394983

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Csf2"
$ws.Range("C5").Value = "Csf2rb"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2300786666666667
$ws.Range("H5").Value = 0.690236
$ws.Range("I5").Value = 0.2137022699341201
$ws.Range("J5").Value = 0.2304051507653011
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 125.4934766666667
$ws.Range("N5").Value = 376.48043
$ws.Range("O5").Value = 0.4861803926831594
$ws.Range("P5").Value = 0.4950178179453587
$ws.Range("Q5").Value = 28.87337178683111
$ws.Range("R5").Value = 259.86034608148
$ws.Range("S5").Value = 0.1038978535138531
$ws.Range("T5").Value = 0.1140546549752107

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Csf2"
$ws.Range("C6").Value = "Csf2rb"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2300786666666667
$ws.Range("H6").Value = 0.690236
$ws.Range("I6").Value = 0.2137022699341201
$ws.Range("J6").Value = 0.2304051507653011
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 36.10423266666667
$ws.Range("N6").Value = 108.312698
$ws.Range("O6").Value = 0.1398731669696947
$ws.Range("P6").Value = 0.1424156772763318
$ws.Range("Q6").Value = 8.306813712969777
$ws.Range("R6").Value = 74.76132341672799
$ws.Range("S6").Value = 0.02989121328429796
$ws.Range("T6").Value = 0.03281330559419569

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Csf2"
$ws.Range("C7").Value = "Csf2rb"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.2341465
$ws.Range("H7").Value = 0.468293
$ws.Range("I7").Value = 0.2174805655477089
$ws.Range("J7").Value = 0.1563191709318771
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.595825
$ws.Range("N7").Value = 27.19165
$ws.Range("O7").Value = 0.05267224809548414
$ws.Range("P7").Value = 0.03575312334118911
$ws.Range("Q7").Value = 3.1834148383625
$ws.Range("R7").Value = 12.73365935345
$ws.Range("S7").Value = 0.01145519030447512
$ws.Range("T7").Value = 0.005588898598919827

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Csf2"
$ws.Range("C8").Value = "Csf2rb"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.2341465
$ws.Range("H8").Value = 0.468293
$ws.Range("I8").Value = 0.2174805655477089
$ws.Range("J8").Value = 0.1563191709318771
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 82.69899733333334
$ws.Range("N8").Value = 248.096992
$ws.Range("O8").Value = 0.3203882151167078
$ws.Range("P8").Value = 0.3262119935919301
$ws.Range("Q8").Value = 19.36368077910933
$ws.Range("R8").Value = 116.182084674656
$ws.Range("S8").Value = 0.06967821021840263
$ws.Range("T8").Value = 0.05099318838632534

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Csf2"
$ws.Range("C9").Value = "Csf2rb"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.2341465
$ws.Range("H9").Value = 0.468293
$ws.Range("I9").Value = 0.2174805655477089
$ws.Range("J9").Value = 0.1563191709318771
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.2286895
$ws.Range("N9").Value = 0.457379
$ws.Range("O9").Value = 0.0008859771349537243
$ws.Range("P9").Value = 0.0006013878451903335
$ws.Range("Q9").Value = 0.05354684601175
$ws.Range("R9").Value = 0.214187384047
$ws.Range("S9").Value = 0.0001926828083720748
$ws.Range("T9").Value = 0.00009400844936866101

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Csf2"
$ws.Range("C10").Value = "Csf2rb"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.2341465
$ws.Range("H10").Value = 0.468293
$ws.Range("I10").Value = 0.2174805655477089
$ws.Range("J10").Value = 0.1563191709318771
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 125.4934766666667
$ws.Range("N10").Value = 376.48043
$ws.Range("O10").Value = 0.4861803926831594
$ws.Range("P10").Value = 0.4950178179453587
$ws.Range("Q10").Value = 29.38385833433167
$ws.Range("R10").Value = 176.30315000599
$ws.Range("S10").Value = 0.1057347867589407
$ws.Range("T10").Value = 0.07738077489772537

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Csf2"
$ws.Range("C11").Value = "Csf2rb"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.2341465
$ws.Range("H11").Value = 0.468293
$ws.Range("I11").Value = 0.2174805655477089
$ws.Range("J11").Value = 0.1563191709318771
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 36.10423266666667
$ws.Range("N11").Value = 108.312698
$ws.Range("O11").Value = 0.1398731669696947
$ws.Range("P11").Value = 0.1424156772763318
$ws.Range("Q11").Value = 8.453679714085668
$ws.Range("R11").Value = 50.722078284514
$ws.Range("S11").Value = 0.03041969545751832
$ws.Range("T11").Value = 0.02226230059953796

# Row 12
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Csf2"
$ws.Range("C12").Value = "Csf2rb"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.6124066666666667
$ws.Range("H12").Value = 1.83722
$ws.Range("I12").Value = 0.5688171645181709
$ws.Range("J12").Value = 0.6132756783028217
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 13.595825
$ws.Range("N12").Value = 27.19165
$ws.Range("O12").Value = 0.05267224809548414
$ws.Range("P12").Value = 0.03575312334118911
$ws.Range("Q12").Value = 8.326173868833335
$ws.Range("R12").Value = 49.957043213
$ws.Range("S12").Value = 0.02996087881047092
$ws.Range("T12").Value = 0.0219265209685122

# Row 13
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Csf2"
$ws.Range("C13").Value = "Csf2rb"
$ws.Range("D13").Value = "Inflammatory-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.6124066666666667
$ws.Range("H13").Value = 1.83722
$ws.Range("I13").Value = 0.5688171645181709
$ws.Range("J13").Value = 0.6132756783028217
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 82.69899733333334
$ws.Range("N13").Value = 248.096992
$ws.Range("O13").Value = 0.3203882151167078
$ws.Range("P13").Value = 0.3262119935919301
$ws.Range("Q13").Value = 50.64541729358223
$ws.Range("R13").Value = 455.80875564224
$ws.Range("S13").Value = 0.1822423160677235
$ws.Range("T13").Value = 0.2000578816406067

# Row 14
$ws.Range("A14").Value = "Neutrophils"
$ws.Range("B14").Value = "Csf2"
$ws.Range("C14").Value = "Csf2rb"
$ws.Range("D14").Value = "MuSCs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6124066666666667
$ws.Range("H14").Value = 1.83722
$ws.Range("I14").Value = 0.5688171645181709
$ws.Range("J14").Value = 0.6132756783028217
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.5
$ws.Range("M14").Value = 0.2286895
$ws.Range("N14").Value = 0.457379
$ws.Range("O14").Value = 0.0008859771349537243
$ws.Range("P14").Value = 0.0006013878451903335
$ws.Range("Q14").Value = 0.1400509743966667
$ws.Range("R14").Value = 0.8403058463799999
$ws.Range("S14").Value = 0.0005039590017323103
$ws.Range("T14").Value = 0.0003688165386821741

# Row 15
$ws.Range("A15").Value = "Neutrophils"
$ws.Range("B15").Value = "Csf2"
$ws.Range("C15").Value = "Csf2rb"
$ws.Range("D15").Value = "Neutrophils"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6124066666666667
$ws.Range("H15").Value = 1.83722
$ws.Range("I15").Value = 0.5688171645181709
$ws.Range("J15").Value = 0.6132756783028217
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 125.4934766666667
$ws.Range("N15").Value = 376.48043
$ws.Range("O15").Value = 0.4861803926831594
$ws.Range("P15").Value = 0.4950178179453587
$ws.Range("Q15").Value = 76.85304173384445
$ws.Range("R15").Value = 691.6773756046
$ws.Range("S15").Value = 0.2765477524103656
$ws.Range("T15").Value = 0.3035823880724225

# Row 16
$ws.Range("A16").Value = "Neutrophils"
$ws.Range("B16").Value = "Csf2"
$ws.Range("C16").Value = "Csf2rb"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6124066666666667
$ws.Range("H16").Value = 1.83722
$ws.Range("I16").Value = 0.5688171645181709
$ws.Range("J16").Value = 0.6132756783028217
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 36.10423266666667
$ws.Range("N16").Value = 108.312698
$ws.Range("O16").Value = 0.1398731669696947
$ws.Range("P16").Value = 0.1424156772763318
$ws.Range("Q16").Value = 22.11047277995111
$ws.Range("R16").Value = 198.99425501956
$ws.Range("S16").Value = 0.07956225822787842
$ws.Range("T16").Value = 0.08734007108259813
